$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the style of the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Populate the time_taken column for each data row (no special style, like the rest of column F/B/C/D data cells)
$ws.Range("F2").Value = "2021-10-05 10:52:40.776815"
$ws.Range("F3").Value = "2021-10-05 10:52:40.776827"
$ws.Range("F4").Value = "2021-10-05 10:52:40.776831"
$ws.Range("F5").Value = "2021-10-05 10:52:40.776835"
$ws.Range("F6").Value = "2021-10-05 10:52:40.776838"
$ws.Range("F7").Value = "2021-10-05 10:52:40.776842"
$ws.Range("F8").Value = "2021-10-05 10:52:40.776845"
$ws.Range("F9").Value = "2021-10-05 10:52:40.776848"
